$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 5838.1934
$ws.Range("I112").Value = 1550
$ws.Range("J112").Value = 5981.1333
$ws.Range("K112").Value = 4650
$ws.Range("L112").Value = 17943.3999
$ws.Range("M112").Value = -3542
$ws.Range("N112").Value = -20159.3999

$ws.Range("H125").Value = 1082.8422
$ws.Range("I125").Value = 1320
$ws.Range("J125").Value = 998.1429000000001
$ws.Range("K125").Value = 11880
$ws.Range("L125").Value = 8983.286100000001
$ws.Range("M125").Value = -9420
$ws.Range("N125").Value = -13903.2861

$ws.Range("H129").Value = 1270.4474
$ws.Range("I129").Value = 556.7
$ws.Range("J129").Value = 1525.3572
$ws.Range("K129").Value = 1670.1
$ws.Range("L129").Value = 4576.071599999999
$ws.Range("M129").Value = 3329.9
$ws.Range("N129").Value = -14576.0716

$ws.Range("H135").Value = 6453.2
$ws.Range("I135").Value = 4218.857
$ws.Range("J135").Value = 11666.667
$ws.Range("K135").Value = 37969.713
$ws.Range("L135").Value = 105000.003
$ws.Range("M135").Value = -35434.713
$ws.Range("N135").Value = -110070.003

$ws.Range("H136").Value = 30000
$ws.Range("J136").Value = 30000
$ws.Range("L136").Value = 30000
$ws.Range("N136").Value = -40200

$ws.Range("H137").Value = 3422.4055
$ws.Range("I137").Value = 2401.3794
$ws.Range("J137").Value = 7123.625
$ws.Range("K137").Value = 7204.138199999999
$ws.Range("L137").Value = 21370.875
$ws.Range("M137").Value = -4654.138199999999
$ws.Range("N137").Value = -26470.875

$ws.Range("H138").Value = 1953.45
$ws.Range("I138").Value = 626.13336
$ws.Range("J138").Value = 2187.6824
$ws.Range("K138").Value = 1878.40008
$ws.Range("L138").Value = 6563.047200000001
$ws.Range("M138").Value = 3261.59992
$ws.Range("N138").Value = -16843.0472

$ws.Range("H141").Value = 5792.7
$ws.Range("I141").Value = 2797.7856
$ws.Range("J141").Value = 12780.833
$ws.Range("K141").Value = 8393.356800000001
$ws.Range("L141").Value = 38342.499
$ws.Range("M141").Value = -3213.356800000001
$ws.Range("N141").Value = -48702.499

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2441.4666
$ws.Range("I61").Value = 1830.7142
$ws.Range("J61").Value = 3866.5557
$ws.Range("K61").Value = 1830.7142
$ws.Range("L61").Value = 3866.5557
$ws.Range("M61").Value = -1618.7142
$ws.Range("N61").Value = -4290.5557

$ws.Range("H74").Value = 2775.0312
$ws.Range("I74").Value = 2516.1765
$ws.Range("J74").Value = 3068.4
$ws.Range("K74").Value = 2516.1765
$ws.Range("L74").Value = 3068.4
$ws.Range("M74").Value = -1642.1765
$ws.Range("N74").Value = -4816.4

$ws.Range("H77").Value = 2775.0312
$ws.Range("I77").Value = 2516.1765
$ws.Range("J77").Value = 3068.4
$ws.Range("K77").Value = 12580.8825
$ws.Range("L77").Value = 15342
$ws.Range("M77").Value = -8212.8825
$ws.Range("N77").Value = -24078

$ws.Range("H110").Value = 1195.7368
$ws.Range("I110").Value = 1224.6471
$ws.Range("J110").Value = 950
$ws.Range("K110").Value = 1224.6471
$ws.Range("L110").Value = 950
$ws.Range("M110").Value = 820.3529000000001
$ws.Range("N110").Value = -5040

$ws.Range("H136").Value = 2441.4666
$ws.Range("I136").Value = 1830.7142
$ws.Range("J136").Value = 3866.5557
$ws.Range("K136").Value = 5492.142599999999
$ws.Range("L136").Value = 11599.6671
$ws.Range("M136").Value = -2942.142599999999
$ws.Range("N136").Value = -16699.6671

$ws.Range("H139").Value = 72598.336
$ws.Range("I139").Value = 40000
$ws.Range("J139").Value = 88897.5
$ws.Range("K139").Value = 40000
$ws.Range("L139").Value = 88897.5
$ws.Range("M139").Value = -34860
$ws.Range("N139").Value = -99177.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2557.147
$ws.Range("I134").Value = 2505.05
$ws.Range("K134").Value = 7515.150000000001
$ws.Range("M134").Value = -4980.150000000001

$ws.Range("H138").Value = 50756.668
$ws.Range("J138").Value = 50756.668
$ws.Range("L138").Value = 50756.668
$ws.Range("N138").Value = -61036.668

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 8157.647
$ws.Range("I31").Value = 2157.6843
$ws.Range("J31").Value = 15757.6
$ws.Range("K31").Value = 2157.6843
$ws.Range("L31").Value = 15757.6
$ws.Range("M31").Value = -1862.6843
$ws.Range("N31").Value = -16347.6

$ws.Range("H34").Value = 8157.647
$ws.Range("I34").Value = 2157.6843
$ws.Range("J34").Value = 15757.6
$ws.Range("K34").Value = 2157.6843
$ws.Range("L34").Value = 15757.6
$ws.Range("M34").Value = -1955.6843
$ws.Range("N34").Value = -16161.6

$ws.Range("H58").Value = 1823.1666
$ws.Range("I58").Value = 1733.8572
$ws.Range("J58").Value = 1880
$ws.Range("K58").Value = 1733.8572
$ws.Range("L58").Value = 1880
$ws.Range("M58").Value = -1530.8572
$ws.Range("N58").Value = -2286

$ws.Range("H132").Value = 9806833
$ws.Range("I132").Value = 2926.5
$ws.Range("J132").Value = 18521418
$ws.Range("K132").Value = 8779.5
$ws.Range("L132").Value = 55564254
$ws.Range("M132").Value = -6249.5
$ws.Range("N132").Value = -55569314

$ws.Range("H134").Value = 2137
$ws.Range("I134").Value = 2126.5454
$ws.Range("J134").Value = 2160
$ws.Range("K134").Value = 6379.6362
$ws.Range("L134").Value = 6480
$ws.Range("M134").Value = -3844.6362
$ws.Range("N134").Value = -11550

$ws.Range("H136").Value = 1823.1666
$ws.Range("I136").Value = 1733.8572
$ws.Range("J136").Value = 1880
$ws.Range("K136").Value = 5201.571599999999
$ws.Range("L136").Value = 5640
$ws.Range("M136").Value = -2651.571599999999
$ws.Range("N136").Value = -10740

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 1012.7778
$ws.Range("I7").Value = 1156.6666
$ws.Range("J7").Value = 725
$ws.Range("K7").Value = 3469.9998
$ws.Range("L7").Value = 2175
$ws.Range("M7").Value = -3357.9998
$ws.Range("N7").Value = -2399

$ws.Range("H34").Value = 11111582
$ws.Range("I34").Value = 97
$ws.Range("J34").Value = 12195629
$ws.Range("K34").Value = 291
$ws.Range("L34").Value = 36586887
$ws.Range("M34").Value = -207
$ws.Range("N34").Value = -36587055

$ws.Range("H39").Value = 1572.2639
$ws.Range("I39").Value = 500
$ws.Range("J39").Value = 1587.3662
$ws.Range("K39").Value = 1500
$ws.Range("L39").Value = 4762.098599999999
$ws.Range("M39").Value = -1206
$ws.Range("N39").Value = -5350.098599999999

$ws.Range("H80").Value = 4400
$ws.Range("J80").Value = 4400
$ws.Range("L80").Value = 13200
$ws.Range("N80").Value = -15072

$ws.Range("H83").Value = 4400
$ws.Range("J83").Value = 4400
$ws.Range("L83").Value = 39600
$ws.Range("N83").Value = -48960

$ws.Range("H107").Value = 553.8421
$ws.Range("I107").Value = 368.83334
$ws.Range("J107").Value = 639.2308
$ws.Range("K107").Value = 1106.50002
$ws.Range("L107").Value = 1917.6924
$ws.Range("M107").Value = 813.4999800000001
$ws.Range("N107").Value = -5757.6924

$ws.Range("H132").Value = 2090.889
$ws.Range("I132").Value = 2575
$ws.Range("J132").Value = 1952.5714
$ws.Range("K132").Value = 23175
$ws.Range("L132").Value = 17573.1426
$ws.Range("M132").Value = -20645
$ws.Range("N132").Value = -22633.1426

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2009.0625
$ws.Range("I126").Value = 2012
$ws.Range("J126").Value = 2002.6
$ws.Range("K126").Value = 6036
$ws.Range("L126").Value = 6007.799999999999
$ws.Range("M126").Value = -3566
$ws.Range("N126").Value = -10947.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 36000
$ws.Range("I93").Value = 100000
$ws.Range("K93").Value = 100000
$ws.Range("M93").Value = -98752

$ws.Range("H132").Value = 3500.6191
$ws.Range("I132").Value = 2713.875
$ws.Range("J132").Value = 3984.7693
$ws.Range("K132").Value = 8141.625
$ws.Range("L132").Value = 11954.3079
$ws.Range("M132").Value = -5611.625
$ws.Range("N132").Value = -17014.3079

$ws.Range("H136").Value = 11113822
$ws.Range("I136").Value = 3075
$ws.Range("J136").Value = 15154094
$ws.Range("K136").Value = 9225
$ws.Range("L136").Value = 45462282
$ws.Range("M136").Value = -6675
$ws.Range("N136").Value = -45467382

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 5211066.5
$ws.Range("I132").Value = 3846.2727
$ws.Range("J132").Value = 7938658
$ws.Range("K132").Value = 11538.8181
$ws.Range("L132").Value = 23815974
$ws.Range("M132").Value = -9008.8181
$ws.Range("N132").Value = -23821034

$ws.Range("H136").Value = 4206.278
$ws.Range("I136").Value = 3977.923
$ws.Range("J136").Value = 4800
$ws.Range("K136").Value = 11933.769
$ws.Range("L136").Value = 14400
$ws.Range("M136").Value = -9383.769
$ws.Range("N136").Value = -19500

Write-Host "Updated Anima_Profits leve sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR)"
